$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Changed skill id to be an enum instead of string": the ##type row (row 3)
# for the id column (B) and dependencies column (E) switches from the
# generic "string" (list) type to the new "SkillId" enum type, and the
# sample data row (row 5) for the ChainLightning skill is rewritten to use
# the new enum-style values.
#
# Touch B5 before B3/E3/E5 so new shared-string entries land in the same
# order Excel produced them in (CHAIN_LIGHTNING, then SkillId, then the
# list type string, then the EXPLOSION;SLOW value).
$ws.Range("B5").Value = "CHAIN_LIGHTNING"
$ws.Range("B3").Value = "SkillId"
$ws.Range("E3").Value = "(list#sep=;),SkillId"
$ws.Range("E5").Value = "EXPLOSION;SLOW"

# Cosmetic formatting nudges that came along with the edit.
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("E5").Font.ThemeColor = 1

# Reflect the cell the editor ended up leaving selected.
$ws.Range("D7").Select()
